$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 459.42856
$ws.Range("I8").Value = 240.76923
$ws.Range("J8").Value = 3302
$ws.Range("K8").Value = 722.30769
$ws.Range("L8").Value = 9906
$ws.Range("M8").Value = -10184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3909.5715
$ws.Range("I64").Value = 3895.05
$ws.Range("J64").Value = 4200
$ws.Range("K64").Value = 3895.05
$ws.Range("L64").Value = 4200
$ws.Range("M64").Value = -3647.05
$ws.Range("N64").Value = -4696

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3909.5715
$ws.Range("I67").Value = 3895.05
$ws.Range("J67").Value = 4200
$ws.Range("K67").Value = 3895.05
$ws.Range("L67").Value = 4200
$ws.Range("M67").Value = -3037.05
$ws.Range("N67").Value = -5916

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2803.25
$ws.Range("I70").Value = 3278.2
$ws.Range("J70").Value = 2464
$ws.Range("K70").Value = 9834.599999999999
$ws.Range("L70").Value = 7392
$ws.Range("M70").Value = -9564.599999999999
$ws.Range("N70").Value = -7932

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2803.25
$ws.Range("I73").Value = 3278.2
$ws.Range("J73").Value = 2464
$ws.Range("K73").Value = 9834.599999999999
$ws.Range("L73").Value = 7392
$ws.Range("M73").Value = -8898.599999999999
$ws.Range("N73").Value = -9264

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1659.5834
$ws.Range("I137").Value = 1464.579
$ws.Range("J137").Value = 2400.6
$ws.Range("K137").Value = 4393.737
$ws.Range("L137").Value = 7201.799999999999
$ws.Range("M137").Value = -1843.737
$ws.Range("N137").Value = -12301.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6275.28
$ws.Range("I61").Value = 7341.1577
$ws.Range("K61").Value = 7341.1577
$ws.Range("M61").Value = -7129.1577

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 111113016
$ws.Range("I63").Value = 125002020
$ws.Range("J63").Value = 1000
$ws.Range("K63").Value = 125002020
$ws.Range("L63").Value = 1000
$ws.Range("M63").Value = -125001334
$ws.Range("N63").Value = -2372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 111113016
$ws.Range("I66").Value = 125002020
$ws.Range("J66").Value = 1000
$ws.Range("K66").Value = 625010100
$ws.Range("L66").Value = 5000
$ws.Range("M66").Value = -625006668
$ws.Range("N66").Value = -11864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1484.6923
$ws.Range("I74").Value = 1476.174
$ws.Range("J74").Value = 1550
$ws.Range("K74").Value = 1476.174
$ws.Range("L74").Value = 1550
$ws.Range("M74").Value = -602.174
$ws.Range("N74").Value = -3298

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1484.6923
$ws.Range("I77").Value = 1476.174
$ws.Range("J77").Value = 1550
$ws.Range("K77").Value = 7380.87
$ws.Range("L77").Value = 7750
$ws.Range("M77").Value = -3012.87
$ws.Range("N77").Value = -16486

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1255.25
$ws.Range("I97").Value = 755
$ws.Range("J97").Value = 1755.5
$ws.Range("K97").Value = 755
$ws.Range("L97").Value = 1755.5
$ws.Range("M97").Value = -259
$ws.Range("N97").Value = -2747.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3082.35
$ws.Range("I132").Value = 1466.9131
$ws.Range("J132").Value = 5267.9414
$ws.Range("K132").Value = 4400.7393
$ws.Range("L132").Value = 15803.8242
$ws.Range("M132").Value = -1870.7393
$ws.Range("N132").Value = -20863.8242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6275.28
$ws.Range("I136").Value = 7341.1577
$ws.Range("K136").Value = 22023.4731
$ws.Range("M136").Value = -19473.4731

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3963.12
$ws.Range("I134").Value = 4746.636
$ws.Range("J134").Value = 2442.1765
$ws.Range("K134").Value = 14239.908
$ws.Range("L134").Value = 7326.529500000001
$ws.Range("M134").Value = -11704.908
$ws.Range("N134").Value = -12396.5295

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 35500
$ws.Range("J14").Value = 35500
$ws.Range("L14").Value = 35500
$ws.Range("N14").Value = -35840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8583.869000000001
$ws.Range("I31").Value = 1929.909
$ws.Range("J31").Value = 14683.333
$ws.Range("K31").Value = 1929.909
$ws.Range("L31").Value = 14683.333
$ws.Range("M31").Value = -1634.909
$ws.Range("N31").Value = -15273.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8583.869000000001
$ws.Range("I34").Value = 1929.909
$ws.Range("J34").Value = 14683.333
$ws.Range("K34").Value = 1929.909
$ws.Range("L34").Value = 14683.333
$ws.Range("M34").Value = -1727.909
$ws.Range("N34").Value = -15087.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1377.9736
$ws.Range("I58").Value = 1059.2858
$ws.Range("J58").Value = 1771.6471
$ws.Range("K58").Value = 1059.2858
$ws.Range("L58").Value = 1771.6471
$ws.Range("M58").Value = -856.2858000000001
$ws.Range("N58").Value = -2177.6471

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13893454
$ws.Range("I99").Value = 1450
$ws.Range("J99").Value = 25007056
$ws.Range("K99").Value = 1450
$ws.Range("L99").Value = 25007056
$ws.Range("M99").Value = 48
$ws.Range("N99").Value = -25010052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 13893454
$ws.Range("I126").Value = 1450
$ws.Range("J126").Value = 25007056
$ws.Range("K126").Value = 4350
$ws.Range("L126").Value = 75021168
$ws.Range("M126").Value = -1880
$ws.Range("N126").Value = -75026108

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2821.95
$ws.Range("I132").Value = 2548.5881
$ws.Range("J132").Value = 4371
$ws.Range("K132").Value = 7645.7643
$ws.Range("L132").Value = 13113
$ws.Range("M132").Value = -5115.7643
$ws.Range("N132").Value = -18173

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1919.6086
$ws.Range("I134").Value = 1832.4117
$ws.Range("K134").Value = 5497.2351
$ws.Range("M134").Value = -2962.2351

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1377.9736
$ws.Range("I136").Value = 1059.2858
$ws.Range("J136").Value = 1771.6471
$ws.Range("K136").Value = 3177.8574
$ws.Range("L136").Value = 5314.9413
$ws.Range("M136").Value = -627.8574000000003
$ws.Range("N136").Value = -10414.9413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 286634.94
$ws.Range("I5").Value = 646.75
$ws.Range("J5").Value = 462627.7
$ws.Range("K5").Value = 1940.25
$ws.Range("L5").Value = 1387883.1
$ws.Range("M5").Value = -1828.25
$ws.Range("N5").Value = -1388107.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 696.3182
$ws.Range("I122").Value = 577.8125
$ws.Range("J122").Value = 1012.3333
$ws.Range("K122").Value = 5200.3125
$ws.Range("L122").Value = 9110.9997
$ws.Range("M122").Value = -2750.3125
$ws.Range("N122").Value = -14010.9997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 286634.94
$ws.Range("I135").Value = 646.75
$ws.Range("J135").Value = 462627.7
$ws.Range("K135").Value = 5820.75
$ws.Range("L135").Value = 4163649.3
$ws.Range("M135").Value = -3285.75
$ws.Range("N135").Value = -4168719.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4189.967
$ws.Range("I132").Value = 4019.3076
$ws.Range("J132").Value = 4320.4707
$ws.Range("K132").Value = 12057.9228
$ws.Range("L132").Value = 12961.4121
$ws.Range("M132").Value = -9527.9228
$ws.Range("N132").Value = -18021.4121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2939
$ws.Range("I61").Value = 2526.5
$ws.Range("J61").Value = 3599
$ws.Range("K61").Value = 2526.5
$ws.Range("L61").Value = 3599
$ws.Range("M61").Value = -2324.5
$ws.Range("N61").Value = -4003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2939
$ws.Range("I113").Value = 2526.5
$ws.Range("J113").Value = 3599
$ws.Range("K113").Value = 2526.5
$ws.Range("L113").Value = 3599
$ws.Range("M113").Value = -356.5
$ws.Range("N113").Value = -7939

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12042561
$ws.Range("I132").Value = 14450207
$ws.Range("J132").Value = 4332.1665
$ws.Range("K132").Value = 43350621
$ws.Range("L132").Value = 12996.4995
$ws.Range("M132").Value = -43348091
$ws.Range("N132").Value = -18056.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5297.8223
$ws.Range("I136").Value = 5420.25
$ws.Range("J136").Value = 5096.1763
$ws.Range("K136").Value = 16260.75
$ws.Range("L136").Value = 15288.5289
$ws.Range("M136").Value = -13710.75
$ws.Range("N136").Value = -20388.5289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1318.0526
$ws.Range("I113").Value = 1265.9286
$ws.Range("J113").Value = 1464
$ws.Range("K113").Value = 3797.7858
$ws.Range("L113").Value = 4392
$ws.Range("M113").Value = -1627.7858
$ws.Range("N113").Value = -8732

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1667.5416
$ws.Range("I132").Value = 1167.7222
$ws.Range("J132").Value = 3167
$ws.Range("K132").Value = 3503.1666
$ws.Range("L132").Value = 9501
$ws.Range("M132").Value = -973.1665999999996
$ws.Range("N132").Value = -14561

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2506.0605
$ws.Range("I136").Value = 2619.44
$ws.Range("J136").Value = 2151.75
$ws.Range("K136").Value = 7858.32
$ws.Range("L136").Value = 6455.25
$ws.Range("M136").Value = -5308.32
$ws.Range("N136").Value = -11555.25
